$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC08 (row 9) - "Insert vehicle with NULL EngineSize for electric vehicle"
# The expectation was flipped: electric vehicles legitimately have no engine
# size, so the insert should succeed (not fail), and the test result flips
# from Pass to Fail, with an updated comment explaining the new expectation.
$ws.Range("F9").Value = "The insertion should be successful as Electric vehicles does not have engine size."
$ws.Range("H9").Value = "Fail"
$ws.Range("I9").Value = "User will not enter any engine size for electric vehicle so we should allow null value."

# Update the view: scroll so column B is the left-most visible column, and
# move the active selection to E19.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E19").Select()
